# Fruta / hortaliza, semanal
# Insert a new data row at row 4 (pushing the existing rows 4-9 down to 5-10)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4; this shifts old rows 4-9 to 5-10
# while carrying formatting (e.g. the date style on column D) along.
$ws.Rows(4).Insert()

# Populate the newly inserted row 4 with the new record's data.
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = "Macroferia Regional de Talca"
$ws.Range("C4").Value = "Maule"
$ws.Range("D4").Value = 45274
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104004
$ws.Range("J4").Value = "Níspero"
$ws.Range("K4").Value = "Golden Nugget"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("Q4").Value = "$/caja 10 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 10
